# New "preseason finished items" added to the FINAL_ITEMS sheet.
# Source: row 107 through row 141 (35 new rows), columns:
#   A = item id, B = item name, C = "Ornn Passive" marker (only for some rows)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FINAL_ITEMS")

$newItems = @(
    @{ Row = 107; Id = 6333; Name = 'Death''s Dance'; Note = $null },
    @{ Row = 108; Id = 6691; Name = 'Duskblade of Draktharr'; Note = 'Ornn Passive' },
    @{ Row = 109; Id = 6616; Name = 'Staff of Flowing Water'; Note = $null },
    @{ Row = 110; Id = 4401; Name = 'Force of Nature'; Note = $null },
    @{ Row = 111; Id = 6609; Name = 'Chempunk Chainsword'; Note = $null },
    @{ Row = 112; Id = 3011; Name = 'Chemtech Putrifier'; Note = $null },
    @{ Row = 113; Id = 4629; Name = 'Cosmic Drive'; Note = $null },
    @{ Row = 114; Id = 4637; Name = 'Demonic Embrace'; Note = $null },
    @{ Row = 115; Id = 6632; Name = 'Divine Sunderer'; Note = 'Ornn Passive' },
    @{ Row = 116; Id = 6692; Name = 'Eclipse'; Note = 'Ornn Passive' },
    @{ Row = 117; Id = 6656; Name = 'Everfrost'; Note = 'Ornn Passive' },
    @{ Row = 118; Id = 6662; Name = 'Frostfire Gauntlet'; Note = 'Ornn Passive' },
    @{ Row = 119; Id = 6671; Name = 'Galeforce'; Note = 'Ornn Passive' },
    @{ Row = 120; Id = 6630; Name = 'Goredrinker'; Note = 'Ornn Passive' },
    @{ Row = 121; Id = 3152; Name = 'Hextech Rocketbelt'; Note = 'Ornn Passive' },
    @{ Row = 122; Id = 4628; Name = 'Horizon Focus'; Note = $null },
    @{ Row = 123; Id = 6673; Name = 'Immortal Shieldbow'; Note = 'Ornn Passive' },
    @{ Row = 124; Id = 4005; Name = 'Imperial Mandate'; Note = 'Ornn Passive' },
    @{ Row = 125; Id = 6672; Name = 'Kraken Slayer'; Note = 'Ornn Passive' },
    @{ Row = 126; Id = 6653; Name = 'Liandry''s Anguish'; Note = 'Ornn Passive' },
    @{ Row = 127; Id = 6655; Name = 'Luden''s Tempest'; Note = 'Ornn Passive' },
    @{ Row = 128; Id = 3222; Name = 'Mikael''s Blessing'; Note = $null },
    @{ Row = 129; Id = 6617; Name = 'Moonstone Renewer'; Note = 'Ornn Passive' },
    @{ Row = 130; Id = 6675; Name = 'Navori Quickblades'; Note = $null },
    @{ Row = 131; Id = 4636; Name = 'Night Harvester'; Note = 'Ornn Passive' },
    @{ Row = 132; Id = 6693; Name = 'Prowler''s Claw'; Note = 'Ornn Passive' },
    @{ Row = 133; Id = 4633; Name = 'Riftmaker'; Note = 'Ornn Passive' },
    @{ Row = 134; Id = 6695; Name = 'Serpent''s Fang'; Note = $null },
    @{ Row = 135; Id = 6694; Name = 'Serylda''s Grudge'; Note = $null },
    @{ Row = 136; Id = 2065; Name = 'Shurelya''s Battlesong'; Note = 'Ornn Passive' },
    @{ Row = 137; Id = 6035; Name = 'Silvermere Dawn'; Note = $null },
    @{ Row = 138; Id = 6631; Name = 'Stridebreaker'; Note = 'Ornn Passive' },
    @{ Row = 139; Id = 3068; Name = 'Sunfire Aegis'; Note = 'Ornn Passive' },
    @{ Row = 140; Id = 6676; Name = 'The Collector'; Note = $null },
    @{ Row = 141; Id = 6664; Name = 'Turbo Chemtank'; Note = 'Ornn Passive' }

)

foreach ($item in $newItems) {
    $ws.Cells.Item($item.Row, 1).Value = $item.Id
    $ws.Cells.Item($item.Row, 2).Value = $item.Name
    if ($item.Note) {
        $ws.Cells.Item($item.Row, 3).Value = $item.Note
    }
}

# Reflect the final on-screen selection from the edit (scrolled down to the
# newly appended rows, with cell G117 active).
$ws.Activate()
$ws.Range("G117").Select()
